# Edit: Thu, Mar 19, 2020  9:05:37 PM
#
# 1) Re-style the comparison table on slide 16 to use the built-in
#    "No Style, No Grid" table style instead of the custom Table_0 style.
# 2) Re-colour the deck's theme: swap the "Integral" colour scheme that
#    currently drives the slide master for the stock "Office" colour
#    scheme (Design > Variants > Colors > "Office").

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------
$tableSlide = $p.Slides.Item(16)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{FF3D5E9A-9970-492C-A3B4-29C508382314}")
    }
}

# --- 2. Theme colours: Integral -> Office ----------------------------------
$themeColors = $p.Slides.Item(1).ThemeColorScheme

# index : scheme slot : target "Office" RGB (stored as 0xBBGGRR, matching
# the VBA RGB()/RGBColor convention used by PowerPoint's ColorFormat.RGB)
$themeColors.Colors(1).RGB  = 0x000000   # dk1
$themeColors.Colors(2).RGB  = 0xFFFFFF   # lt1
$themeColors.Colors(3).RGB  = 0x6A5444   # dk2
$themeColors.Colors(4).RGB  = 0xE6E6E7   # lt2
$themeColors.Colors(5).RGB  = 0xD59B5B   # accent1
$themeColors.Colors(6).RGB  = 0x317DED   # accent2
$themeColors.Colors(7).RGB  = 0xA5A5A5   # accent3
$themeColors.Colors(8).RGB  = 0x00C0FF   # accent4
$themeColors.Colors(9).RGB  = 0xC47244   # accent5
$themeColors.Colors(10).RGB = 0x47AD70   # accent6
$themeColors.Colors(11).RGB = 0xC16305   # hlink
$themeColors.Colors(12).RGB = 0x724F95   # folHlink
